# New watchlist script implementation
# Adds a new test case row (Watchlist035 / OPQA-620) to the "Test Cases" sheet,
# mirroring the existing rows' layout, styling and shared-string usage, and
# updates the sheet's view/selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 36
$prevRow = 35

# --- Copy the formatting of the previous row's cells onto the new row ---
# Columns A, B, C and E follow row 35's styling; column D follows row 34's
# styling (it carries a slightly different fill/border combination).
$ws.Range("A" + $prevRow).Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)

$ws.Range("B" + $prevRow).Copy()
$ws.Range("B" + $newRow).PasteSpecial(-4122)

$ws.Range("C" + $prevRow).Copy()
$ws.Range("C" + $newRow).PasteSpecial(-4122)

$ws.Range("D" + ($prevRow - 1)).Copy()
$ws.Range("D" + $newRow).PasteSpecial(-4122)

$ws.Range("E" + $prevRow).Copy()
$ws.Range("E" + $newRow).PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# --- Fill in the new test case's values ---
# Values are written in Jira id -> Description -> TCID order so the shared
# string table grows in the same sequence as the source workbook.
$ws.Range("B" + $newRow).Value2 = "OPQA-620"
$ws.Range("C" + $newRow).Value2 = "Verify that user is able to comment on an item contained in public watchlist of some other user"
$ws.Range("A" + $newRow).Value2 = "Watchlist035"
$ws.Range("D" + $newRow).Value2 = "Y"

# --- Update the sheet view: scroll position and active selection ---
$excel.Goto($ws.Range("A31"), $true) | Out-Null
$aw = $excel.ActiveWindow
$aw.ScrollRow = 31
$aw.ScrollColumn = 1
$ws.Range("C38").Select() | Out-Null
